$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New empatica "file splitter" test recordings for adri (adriana) and leile (leichtle),
# appended below the existing roster (rows 2-13) starting at row 14.
$newRows = @(
    @("adriana",  "A02F50"),
    @("adriana4", "A02F50"),
    @("adriana2", "A02F50"),
    @("adriana1", "A02F50"),
    @("leichtle",  "A02F6F"),
    @("leichtle4", "A02F6F"),
    @("leichtle2", "A02F6F"),
    @("leichtle1", "A02F6F")
)

$r = 14
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Widen the empa-id column (C) to fit the new longer device ids.
$ws.Columns.Item(3).ColumnWidth = 14.16

# Leave row 23 touched (blank spacer row below the new data) with its own row height.
$ws.Range("A23").NumberFormat = "General"
$ws.Rows.Item(23).RowHeight = 13.8

$ws.Range("C27").Select()
